$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Fill in the missing column B values for rows 1-5, matching the pattern
# already present from row 6 onward (B = "karbohidrat"), except B1 which
# mirrors A1's own text ("q").
$ws.Range("B1").Value = "q"
$ws.Range("B2").Value = "karbohidrat"
$ws.Range("B3").Value = "karbohidrat"
$ws.Range("B4").Value = "karbohidrat"
$ws.Range("B5").Value = "karbohidrat"

# Update the saved selection/view: move the selection to B2.
$ws.Activate()
$ws.Range("B2").Select()
